$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target B:E values for rows 2-51 (row index, B, C, D, E)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '43.123.56', '  +2.35%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '2.305.76', '  +1.66%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.02%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '301.59', '  +1.19%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '99.33', '  +5.44%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.506', '  +2.42%  '),
    @(8, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  -0.05%  '),
    @(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.505', '  +3.08%  '),
    @(10, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '34.19', '  +3.69%  '),
    @(11, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.0799', '  +1.22%  '),
    @(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.117', '  +4.13%  '),
    @(13, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '17.92', '  +14.94%  '),
    @(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '6.81', '  +2.35%  '),
    @(15, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.663.58', '  +1.71%  '),
    @(16, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '2.279.70', '  +1.04%  '),
    @(17, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.811', '  +4.99%  '),
    @(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '43.031.36', '  +2.14%  '),
    @(19, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '12.56', '  +10.74%  '),
    @(20, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0₃0908', '  +2.06%  '),
    @(21, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.09', '  +1.88%  '),
    @(22, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '67.76', '  +1.65%  '),
    @(23, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '236.79', '  +1.47%  '),
    @(24, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.21', '  +14.77%  '),
    @(25, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.00', '  +0.08%  '),
    @(26, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '2.46', '  +0.64%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '24.78', '  +4.19%  '),
    @(28, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '168.15', '  +0.71%  '),
    @(29, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.09', '  -9.13%  '),
    @(30, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '34.07', '  +1.32%  '),
    @(31, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.14', '  +1.31%  '),
    @(32, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '0.999', '  -0.01%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.06', '  +2.86%  '),
    @(34, 'WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '2.44', '  +4.34%  '),
    @(35, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '4.55', '  +1.69%  '),
    @(36, 'Celestia', 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia', '16.97', '  +5.92%  '),
    @(37, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0691', '  +0.19%  '),
    @(38, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.103', '  +3.88%  '),
    @(39, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.80', '  +5.49%  '),
    @(40, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.81', '  +1.19%  '),
    @(41, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.110', '  +0.59%  '),
    @(42, 'ApeXProtocol', 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex', '2.35', '  -3.83%  '),
    @(43, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.994.76', '  +1.67%  '),
    @(44, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0285', '  +2.89%  '),
    @(45, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '10.06', '  +5.49%  '),
    @(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '17.76', '  +3.21%  '),
    @(47, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '2.85', '  +2.71%  '),
    @(48, 'MultiversX', 'https://coinranking.com/coin/omwkOTglq+multiversx-egld', '56.56', '  +9.07%  '),
    @(49, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.530.70', '  +1.54%  '),
    @(50, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '1.54', '  +4.91%  '),
    @(51, 'THORChain', 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune', '4.56', '  +1.04%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    # Column D: prefix with apostrophe to force text storage (avoid numeric auto-conversion),
    # then reset style so no stray NumberFormat/quotePrefix style is left on the cell.
    $ws.Cells.Item($r, 4).Value = "'" + $row[3]
    $ws.Cells.Item($r, 4).Style = 'Normal'
    $ws.Cells.Item($r, 5).Value = $row[4]
}